$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values of the columns that change, for rows 2-7.
# Note: use Value2 (not Value) to reliably read back plain scalars in this runtime.
$cols = @("D","L","M","N","O","P","R","S")
$orig = @{}
foreach ($r in 2..7) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $row
}

# Row permutation: new row N gets the data that used to live in old row Map[N]
$map = @{ 2 = 6; 3 = 7; 4 = 5; 5 = 4; 6 = 2; 7 = 3 }

foreach ($newRow in 2..7) {
    $srcRow = $map[$newRow]
    $src = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
